$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: merge the two bold runs "Weighted Average " + "(user manual P5)"
# into a single run "Weighted Average (user manual P5)".
# ------------------------------------------------------------------
$d.Content.Find.Execute("Weighted Average (user manual P5)", $false, $false, $false, $false, $false, $true, 1, $false, "Weighted Average (user manual P5)", 2) | Out-Null

# ------------------------------------------------------------------
# Change 2: rewrite the three paragraphs following "factor score to
# quintiles" (that paragraph itself, plus the two blank red paragraphs
# after it) into two paragraphs:
#   - a paragraph with the full explanatory sentence about
#     factorScore_to_quintiles.py
#   - a paragraph with three runs: a dashed divider, "Kenny reflect to
#     here", and another dashed divider
# ------------------------------------------------------------------
$pStart = $d.Paragraphs(38)
$pEnd = $d.Paragraphs(40)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$xmlFrag = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Based on user manual P5, after calculating weighted average of factor score, use factorScore_to_quintiles.py to get  quintiles in 2016_weightedAverage_calculated_quintiles.csv</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>---------------------------------------------</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Kenny reflect to here</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>------------------------------------</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xmlFrag)

# ------------------------------------------------------------------
# Incidental styles.xml change: a "ListLabel 6" character style shows up
# in the saved stylesheet (mirroring ListLabel 1..5 already present).
# ------------------------------------------------------------------
$listLabel6 = $d.Styles.Add("ListLabel 6", 2)
$listLabel6.QuickStyle = $true
